$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("testdata_Prop_1")

$ws1.Range("E1").Value = "Proportion"
$ws1.Range("F1").Value = "lower95cl"
$ws1.Range("G1").Value = "upper95cl"
$ws1.Range("H1").Value = "lower998cl"
$ws1.Range("I1").Value = "upper998cl"

$ws1.Activate()
$ws1.Range("F5").Select()
